# Atualização automática da planilha
# Adds 7 new "Time Projeto" / Key User rows (36-42) to the Organograma sheet,
# reusing the existing visual style of the table (green-filled band with thin
# ECEAE6 borders) and alternating between a "full border" row style and a
# "grouped / no top border" row style, matching the source edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Organograma")
$nbsp = [char]0x00A0

function Add-FullRow($row, $name, $area) {
    $ws.Range("C35:E35").Copy() | Out-Null
    $ws.Range("C" + $row).PasteSpecial(-4122) | Out-Null
    $ws.Range("C" + $row).Value2 = $name
    $ws.Range("D" + $row).Value2 = $nbsp
    $ws.Range("E" + $row).Value2 = $area
}

function Add-BandedRow($row, $name, $area) {
    $ws.Range("C35:E35").Copy() | Out-Null
    $ws.Range("C" + $row).PasteSpecial(-4122) | Out-Null
    $ws.Range("C" + $row + ":E" + $row).Borders(8).LineStyle = -4142
    $ws.Range("C" + $row).Value2 = $name
    $ws.Range("D" + $row).Value2 = $nbsp
    $ws.Range("E" + $row).Value2 = $area
}

Add-FullRow   36 "Vanessa Kato"     "Stand"
Add-FullRow   37 "Marcos Souza"     "Viabilidade Economica"
Add-BandedRow 38 "Bruna Fernandes"  "Aprovações"
Add-BandedRow 39 "Fabiana Passos"   "Inteligência de Mercado"
Add-FullRow   40 "Sandra Trombeli"  "Projetos"
Add-BandedRow 41 "Carlos Morais"    "Infraestrutura"
Add-BandedRow 42 "Silvia Reis"      "Custos"

$excel.CutCopyMode = 0
